$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBNPPTY")

# Historical years (2021-2023, columns B:D) should now be treated like the
# "always allowed" boolean flag (value 1) for every technology row, and
# lose the special "0/1 flag" number format (style index 5) that previously
# formatted them - matching the plain/general formatting used elsewhere
# in the model once the flag is fixed at 1.
$rng = $ws.Range("B2:D25")
$rng.Value = 1
$rng.Style = "Normal"

# Reflect the updated selection/active cell left behind on this sheet.
$ws.Range("B2:D25").Select() | Out-Null
